$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.550.71'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '1.849.41'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '262.51'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5252'
$ws.Range("E7").Value = '  +1.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3233'
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06799'
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("E10").Value = '  +1.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7818'
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07767'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("D13").Value = '1.851.06'
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.62'
$ws.Range("E14").Value = '  -0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.031'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.95'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007957'
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("D20").Value = '26.571.01'
$ws.Range("E20").Value = '  +0.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.642'
$ws.Range("E21").Value = '  +2.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.454'
$ws.Range("E22").Value = '  +0.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.997'
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.06'
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.166'
$ws.Range("E25").Value = '  -4.95%  '
$ws.Range("E26").Value = '  +2.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.01'
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '112.33'
$ws.Range("E28").Value = '  +1.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.177'
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08727'
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.100'
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04871'
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.131'
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7196'
$ws.Range("E34").Value = '  +5.49%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.877'
$ws.Range("E35").Value = '  +1.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.103'
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.279'
$ws.Range("E37").Value = '  +2.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01784'
$ws.Range("E38").Value = '  +0.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.4855'
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8995'
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.955'
$ws.Range("E42").Value = '  -2.63%  '
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.680'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4170'
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05869'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.034'
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.16'
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1233'
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.8957'
$ws.Range("E50").Value = '  +3.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.93'
$ws.Range("E51").Value = '  +1.10%  '
